# "230523 Tools Status.xlsx" - update the Zipline Launcher status note and
# the sheet's cached selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# E8 held the placeholder note "not implemented". Give it the same
# "Neutral" wrap-text note style already used by the other multi-line
# status cells (copy F8's formatting), then replace the text with the
# real status note.
$null = $ws.Range("F8").Copy()
$null = $ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "model missing;
rotation missing;"

# The workbook's cached selection moves from E8 to E11.
$null = $ws.Range("E11").Select()
